$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three Neo4j query strings (CasesTab / SamplesTab / FilesTab) ---
# B2 (CasesTab query): append an ORDER BY / LIMIT clause on a new line
$c2 = $ws.Range("B2")
$c2.Value = $c2.Value() + "`n order By ss.study_subject_id ASC LIMIT 100"

# B3 (SamplesTab query): append an ORDER BY / LIMIT clause on a new line
$c3 = $ws.Range("B3")
$c3.Value = $c3.Value() + "`n order By samp.sample_id ASC LIMIT 100"

# B4 (FilesTab query): replace the trailing "    order by f.file_name" with
# "     order By f.file_name ASC LIMIT 100"
$c4 = $ws.Range("B4")
$v4 = $c4.Value()
$oldTail = "    order by f.file_name"
$newTail = "     order By f.file_name ASC LIMIT 100"
$v4 = $v4.Substring(0, $v4.Length - $oldTail.Length) + $newTail
$c4.Value = $v4

# --- Row height changes caused by the extra wrapped line of text ---
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 374.4

# --- Update selection / active cell to B4 ---
$ws.Range("B4").Select() | Out-Null
